# Add an "Address" column to the worksheet.
# The existing column F ("District") is shifted one place to the right,
# to column G, and a new column F ("Address") is populated with
# per-row address values (mostly extracted from the second line of the
# Names column, B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; this pushes the existing F column
# (District) to G and everything after it keeps moving right.
$ws.Range("F:F").Insert()

# Header for the new column.
$ws.Range("F2").Value = "Address"

# Per-row address values.
$addresses = @{
    3  = "Govt. Adarsh Vidyalaya Indi"
    4  = "Govt. High School (RMSA) Jaddigadde"
    5  = "Govt. High School Indi"
    6  = "Smt. Ramabai High School Indi"
    7  = "B E H S NidagundiB Bagewadi"
    8  = "Talikot Anjuman High SchoolMudebihal"
    9  = "G H S AnashiJoida"
    10 = "Govt. High SchoolAlurIndi"
    11 = "Govt. High School AmadalliKarwar"
    12 = "V Y Patil High School PadanurIndi"
    13 = "Govt. High SchoolKangodSiddapur"
    14 = "S Y High School MiragiIndi"
    15 = "Govt. High School BaragudiIndi"
    16 = "Govt. High School ShiralagiSidddapur"
    17 = "G H S PeerapurMuddebhihal"
    18 = "Rahimkhan Unity High School Karwar"
    19 = "G H S (RMSA) ArjanalIndiChadachan"
    20 = "R P S S Secondary School BelseAnkola"
    21 = "Muddebihal"
    22 = "Govt. High School MadikeshwarMuddebihal"
    23 = "The PNE SchoolCharliaKarwar"
    24 = "S K Comp. P U College Talikoti"
    25 = "B J Patil Composite P U CollegeManaguliB Bagewadi"
    26 = "Govt. High School AllankiHonavarKarwar"
    27 = "Govt. High School B SalawadgiMuddebihal"
    29 = "D B E H S Nimbal LtIndi"
    30 = "Janata VidyalayaAnilgodHonnavar"
    31 = "K R C R School TavanchuruSiddapur"
    32 = "G H S HinganiCHadachan"
    34 = "Siddapur"
    35 = "P S C H S Sindagi"
}
# Rows 28 and 33 intentionally have no address value (column F stays
# blank there, matching the source data).

foreach ($row in $addresses.Keys) {
    $ws.Range("F$row").Value = $addresses[$row]
}

$ws.Range("A1").Select()
